# Apply the "address change" edit described in the commit:
# - Street changed from "Neuweilerstrasse 88" to "Im Westfeld 17"
# - Postal code changed from "4054 Basel" to "4055 Basel"
# - A new "MwSt Nr: CHE-108.242.406" line was added below the bank details
# - The "Summe Total:" label (near the grand total) was moved one column to
#   the left (from M42 to L42) and its text updated to mention the 7.7% VAT

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the VAT number line under the bank/account information block
# (copy the formatting of the line above it, then set the new text)
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = "MwSt Nr: CHE-108.242.406"

# --- Move the "Summe Total" label and update its wording for the VAT note
$ws.Range("M42").Copy($ws.Range("L42"))
$ws.Range("M42").Clear()
$ws.Range("L42").Value = "Summe Total inkl. 7.7% MwSt:"

# --- Update the company address (it appears twice: header block + footer block)
$ws.Range("A3").Value = "Im Westfeld 17"
$ws.Range("A16").Value = "Im Westfeld 17"
$ws.Range("A17").Value = "4055 Basel"

# --- Reflect where the author was working in the sheet when saving
$ws.Range("A20").Select()

# --- Page margins were reset to Excel's normal defaults
$ws.PageSetup.LeftMargin = 18
$ws.PageSetup.RightMargin = 18
$ws.PageSetup.TopMargin = 54
$ws.PageSetup.BottomMargin = 54
$ws.PageSetup.HeaderMargin = 21.6
$ws.PageSetup.FooterMargin = 21.6

$wb.Save()
